# Updates cryptocurrency price/volume figures in the tracking sheet
# (commit: 'Updated cryptos list ... with GitHub Actions').
# Column D (Price) values are text that look numeric (e.g. '7.81' or
# '66.735.10'), so a leading apostrophe forces Excel to keep them as
# literal text instead of auto-converting to a Number/Date.
# Column E (Volume 1h) values already contain '%', '+' and padding
# spaces so Excel stores them as text without any extra nudging.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates ---
$ws.Range("D2").Value = "'66.735.10"
$ws.Range("D3").Value = "'3.495.82"
$ws.Range("D6").Value = "'147.38"
$ws.Range("D7").Value = "'3.494.52"
$ws.Range("D11").Value = "'7.81"
$ws.Range("D14").Value = "'4.086.85"
$ws.Range("D15").Value = "'31.25"
$ws.Range("D16").Value = "'3.486.63"
$ws.Range("D17").Value = "'66.757.10"
$ws.Range("D19").Value = "'10.51"
$ws.Range("D21").Value = "'15.33"
$ws.Range("D23").Value = "'0.608"
$ws.Range("D24").Value = "'79.73"
$ws.Range("D25").Value = "'3.632.29"
$ws.Range("D29").Value = "'9.77"
$ws.Range("D30").Value = "'8.22"
$ws.Range("D33").Value = "'1.00"
$ws.Range("D35").Value = "'25.29"
$ws.Range("D36").Value = "'3.489.24"
$ws.Range("D39").Value = "'8.00"
$ws.Range("D41").Value = "'0.999"
$ws.Range("D43").Value = "'169.83"
$ws.Range("D45").Value = "'5.40"
$ws.Range("D47").Value = "'29.00"
$ws.Range("D49").Value = "'1.30"

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("E6").Value = "  -3.22%  "
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  -1.70%  "
$ws.Range("E10").Value = "  -0.99%  "
$ws.Range("E11").Value = "  +2.23%  "
$ws.Range("E12").Value = "  -2.25%  "
$ws.Range("E13").Value = "  -1.49%  "
$ws.Range("E14").Value = "  -0.45%  "
$ws.Range("E15").Value = "  -4.66%  "
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("E17").Value = "  -1.47%  "
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("E19").Value = "  +6.62%  "
$ws.Range("E20").Value = "  -2.81%  "
$ws.Range("E21").Value = "  -1.64%  "
$ws.Range("E22").Value = "  -3.61%  "
$ws.Range("E23").Value = "  -3.94%  "
$ws.Range("E24").Value = "  +1.93%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E28").Value = "  -6.29%  "
$ws.Range("E29").Value = "  -3.20%  "
$ws.Range("E30").Value = "  -7.58%  "
$ws.Range("E31").Value = "  -0.81%  "
$ws.Range("E32").Value = "  -4.63%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  -1.62%  "
$ws.Range("E35").Value = "  -1.89%  "
$ws.Range("E36").Value = "  -0.44%  "
$ws.Range("E37").Value = "  -4.85%  "
$ws.Range("E38").Value = "  -5.33%  "
$ws.Range("E39").Value = "  -0.89%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("E42").Value = "  -1.09%  "
$ws.Range("E43").Value = "  -2.47%  "
$ws.Range("E44").Value = "  -9.71%  "
$ws.Range("E45").Value = "  -1.67%  "
$ws.Range("E46").Value = "  +1.61%  "
$ws.Range("E47").Value = "  -4.51%  "
$ws.Range("E48").Value = "  -2.61%  "
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("E50").Value = "  -3.03%  "
$ws.Range("E51").Value = "  -4.45%  "
